# 25.06.2023 mybee end rewrite
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 8 ----
$ws.Cells.Item(8, 5).Value  = 10        # E8
$ws.Cells.Item(8, 6).Value  = 288.73    # F8
$ws.Cells.Item(8, 8).Value  = 10        # H8
$ws.Cells.Item(8, 9).Value  = 831.58    # I8
$ws.Cells.Item(8, 11).Value = 4         # K8
$ws.Cells.Item(8, 12).Value = 116.73    # L8

# ---- Row 9 ----
$ws.Cells.Item(9, 5).Value  = 90        # E9
$ws.Cells.Item(9, 6).Value  = 2598.59   # F9
$ws.Cells.Item(9, 8).Value  = 90        # H9
$ws.Cells.Item(9, 9).Value  = 7484.25   # I9
$ws.Cells.Item(9, 11).Value = 96        # K9
$ws.Cells.Item(9, 12).Value = 2801.58   # L9
$ws.Cells.Item(9, 16).Value = 1920.11   # P9
$ws.Cells.Item(9, 18).Value = 0.93      # R9
$ws.Cells.Item(9, 20).Value = 1.06      # T9

# ---- Row 11 ----
$ws.Cells.Item(11, 11).Value = 2                  # K11
$ws.Cells.Item(11, 12).Value = 72.81              # L11
$ws.Cells.Item(11, 16).Value = 1935.94            # P11
$ws.Cells.Item(11, 18).Value = 0.9399999999999999 # R11
$ws.Cells.Item(11, 20).Value = 1.06               # T11

# ---- Row 12 ----
$ws.Cells.Item(12, 11).Value = 98       # K12
$ws.Cells.Item(12, 12).Value = 3567.61  # L12

# ---- Row 35 ----
$ws.Cells.Item(35, 16).Value = 22349.13 # P35
$ws.Cells.Item(35, 18).Value = 10.8     # R35
$ws.Cells.Item(35, 20).Value = 12.28    # T35

# ---- Rows 41-47 rewritten, rows 48-49 removed ----
$ws.Cells.Item(41, 1).Value = "Аккумуляторный"
$ws.Cells.Item(41, 2).Value = 21
$ws.Cells.Item(41, 3).Value = 15
$ws.Cells.Item(41, 4).Value = 1
$ws.Cells.Item(41, 5).Value = 21

$ws.Cells.Item(42, 1).Value = "Топливной аппаратуры"
$ws.Cells.Item(42, 2).Value = 14
$ws.Cells.Item(42, 3).Value = 8
$ws.Cells.Item(42, 4).Value = 1
$ws.Cells.Item(42, 5).Value = 14

$ws.Cells.Item(43, 1).Value = "Шиномонтажный"
$ws.Cells.Item(43, 2).Value = 18
$ws.Cells.Item(43, 3).Value = 15
$ws.Cells.Item(43, 4).Value = 1
$ws.Cells.Item(43, 5).Value = 18

$ws.Cells.Item(44, 1).Value = "Кузнечно-рессорный"
$ws.Cells.Item(44, 2).Value = 21
$ws.Cells.Item(44, 3).Value = 5
$ws.Cells.Item(44, 4).Value = 1
$ws.Cells.Item(44, 5).Value = 21

$ws.Cells.Item(45, 1).Value = "Сварочный"
$ws.Cells.Item(45, 2).Value = 15
$ws.Cells.Item(45, 3).Value = 9
$ws.Cells.Item(45, 4).Value = 1
$ws.Cells.Item(45, 5).Value = 15

$ws.Cells.Item(46, 1).Value = "Жестяницкий"
$ws.Cells.Item(46, 2).Value = 18
$ws.Cells.Item(46, 3).Value = 12
$ws.Cells.Item(46, 4).Value = 1
$ws.Cells.Item(46, 5).Value = 18

$ws.Cells.Item(47, 1).Value = "Малярный"
$ws.Cells.Item(47, 2).Value = 30
$ws.Cells.Item(47, 3).Value = 15
$ws.Cells.Item(47, 4).Value = 2
$ws.Cells.Item(47, 5).Value = 45

# Rows 48 and 49 no longer exist in the rewritten table - clear them out
$ws.Range("A48:T49").ClearContents()

# ---- New rows 112-155 ----
$newRows = @(
    @(112, "УАЗ-3163",        29,  2,    0.85, 1.4,  1,   1, 1.1, 8),
    @(113, "КАМАЗ-43502",     124, 4,    0.8,  1.15, 0.8, 1, 1.1, 40),
    @(114, "КАМАЗ-43118",     144, 4,    0.85, 1.15, 1.3, 1, 1.1, 81),
    @(115, "ПОЛИТРАНС-94163", 61,  1,    0.85, 1.2,  1.5, 1, 1.1, 10),
    @(116, "УАЗ-3163",        29,  1.5,  0.85, 1.4,  1,   1, 1.1, 6),
    @(117, "КАМАЗ-43502",     124, 2.5,  0.8,  1.15, 0.8, 1, 1.1, 25),
    @(118, "КАМАЗ-43118",     144, 2.5,  0.85, 1.15, 1.3, 1, 1.1, 50),
    @(119, "ПОЛИТРАНС-94163", 61,  0,    0.85, 1.2,  1.5, 1, 1.1, 1),
    @(120, "УАЗ-3163",        29,  1.5,  0.85, 1.4,  1,   1, 1.1, 6),
    @(121, "КАМАЗ-43502",     124, 1.6,  0.8,  1.15, 0.8, 1, 1.1, 16),
    @(122, "КАМАЗ-43118",     144, 1.6,  0.85, 1.15, 1.3, 1, 1.1, 32),
    @(123, "ПОЛИТРАНС-94163", 61,  0.3,  0.85, 1.2,  1.5, 1, 1.1, 3),
    @(124, "УАЗ-3163",        29,  0.4,  0.85, 1.4,  1,   1, 1.1, 2),
    @(125, "КАМАЗ-43502",     124, 0.5,  0.8,  1.15, 0.8, 1, 1.1, 5),
    @(126, "КАМАЗ-43118",     144, 0.5,  0.85, 1.15, 1.3, 1, 1.1, 10),
    @(127, "ПОЛИТРАНС-94163", 61,  0.2,  0.85, 1.2,  1.5, 1, 1.1, 2),
    @(128, "УАЗ-3163",        29,  0.1,  0.85, 1.4,  1,   1, 1.1, 1),
    @(129, "КАМАЗ-43502",     124, 0.15, 0.8,  1.15, 0.8, 1, 1.1, 2),
    @(130, "КАМАЗ-43118",     144, 0.15, 0.85, 1.15, 1.3, 1, 1.1, 3),
    @(131, "ПОЛИТРАНС-94163", 61,  0.05, 0.85, 1.2,  1.5, 1, 1.1, 1),
    @(132, "УАЗ-3163",        29,  0.15, 0.85, 1.4,  1,   1, 1.1, 1),
    @(133, "КАМАЗ-43502",     124, 0.15, 0.8,  1.15, 0.8, 1, 1.1, 2),
    @(134, "КАМАЗ-43118",     144, 0.15, 0.85, 1.15, 1.3, 1, 1.1, 3),
    @(135, "ПОЛИТРАНС-94163", 61,  0.1,  0.85, 1.2,  1.5, 1, 1.1, 1),
    @(136, "УАЗ-3163",        29,  0,    0.85, 1.4,  1,   1, 1.1, 1),
    @(137, "КАМАЗ-43502",     124, 0.3,  0.8,  1.15, 0.8, 1, 1.1, 3),
    @(138, "КАМАЗ-43118",     144, 0.3,  0.85, 1.15, 1.3, 1, 1.1, 6),
    @(139, "ПОЛИТРАНС-94163", 61,  0.2,  0.85, 1.2,  1.5, 1, 1.1, 2),
    @(140, "УАЗ-3163",        29,  0.2,  0.85, 1.4,  1,   1, 1.1, 1),
    @(141, "КАМАЗ-43502",     124, 0.25, 0.8,  1.15, 0.8, 1, 1.1, 3),
    @(142, "КАМАЗ-43118",     144, 0.25, 0.85, 1.15, 1.3, 1, 1.1, 5),
    @(143, "ПОЛИТРАНС-94163", 61,  0.15, 0.85, 1.2,  1.5, 1, 1.1, 2),
    @(144, "УАЗ-3163",        29,  1.6,  0.85, 1.4,  1,   1, 1.1, 6),
    @(145, "КАМАЗ-43502",     124, 2.4,  0.8,  1.15, 0.8, 1, 1.1, 24),
    @(146, "КАМАЗ-43118",     144, 2.4,  0.85, 1.15, 1.3, 1, 1.1, 48),
    @(147, "ПОЛИТРАНС-94163", 61,  1.2,  0.85, 1.2,  1.5, 1, 1.1, 12),
    @(148, "УАЗ-3163",        29,  4,    0.85, 1.4,  1,   1, 1.1, 15),
    @(149, "КАМАЗ-43502",     124, 6,    0.8,  1.15, 0.8, 1, 1.1, 60),
    @(150, "КАМАЗ-43118",     144, 6,    0.85, 1.15, 1.3, 1, 1.1, 121),
    @(151, "ПОЛИТРАНС-94163", 61,  2,    0.85, 1.2,  1.5, 1, 1.1, 21),
    @(152, "УАЗ-3163",        29,  0.4,  0.85, 1.4,  1,   1, 1.1, 2),
    @(153, "КАМАЗ-43502",     124, 0.8,  0.8,  1.15, 0.8, 1, 1.1, 8),
    @(154, "КАМАЗ-43118",     144, 0.8,  0.85, 1.15, 1.3, 1, 1.1, 16),
    @(155, "ПОЛИТРАНС-94163", 61,  0.2,  0.85, 1.2,  1.5, 1, 1.1, 2)
)

foreach ($r in $newRows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 1).Value = $r[1]
    $ws.Cells.Item($rowNum, 2).Value = $r[2]
    $ws.Cells.Item($rowNum, 3).Value = $r[3]
    $ws.Cells.Item($rowNum, 4).Value = $r[4]
    $ws.Cells.Item($rowNum, 5).Value = $r[5]
    $ws.Cells.Item($rowNum, 6).Value = $r[6]
    $ws.Cells.Item($rowNum, 7).Value = $r[7]
    $ws.Cells.Item($rowNum, 8).Value = $r[8]
    $ws.Cells.Item($rowNum, 9).Value = $r[9]
}
